$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053326784846613
$ws.Range("D2").Value = 1.05735288753316
$ws.Range("E2").Value = 1.049784570808107
$ws.Range("F2").Value = 1.065580750290416
$ws.Range("I2").Value = 1.042826644038349
$ws.Range("J2").Value = 1.05834471400405
$ws.Range("K2").Value = 1.060088090727344
$ws.Range("L2").Value = 1.052540665006014
$ws.Range("M2").Value = 1.068293606303547

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055145191482077
$ws.Range("D3").Value = 1.058813802458066
$ws.Range("E3").Value = 1.051380687837004
$ws.Range("F3").Value = 1.067271853948691
$ws.Range("I3").Value = 1.043339064684927
$ws.Range("J3").Value = 1.059809744085158
$ws.Range("K3").Value = 1.061361181682897
$ws.Range("L3").Value = 1.053947076196142
$ws.Range("M3").Value = 1.069797949821408

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056318172525178
$ws.Range("D4").Value = 1.05975583367735
$ws.Range("E4").Value = 1.052409786857123
$ws.Range("F4").Value = 1.068363033590878
$ws.Range("I4").Value = 1.043667506075709
$ws.Range("J4").Value = 1.060753808169771
$ws.Range("K4").Value = 1.062181170397339
$ws.Range("L4").Value = 1.054852949132864
$ws.Range("M4").Value = 1.070767819387793

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056810438709082
$ws.Range("D5").Value = 1.060151092597325
$ws.Range("E5").Value = 1.052841551764209
$ws.Range("F5").Value = 1.06882104461369
$ws.Range("I5").Value = 1.043804839804519
$ws.Range("J5").Value = 1.061149772308975
$ws.Range("K5").Value = 1.06252500024025
$ws.Range("L5").Value = 1.055232795306993
$ws.Range("M5").Value = 1.071174719554557

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056893042766477
$ws.Range("D6").Value = 1.06021741349697
$ws.Range("E6").Value = 1.052913996545844
$ws.Range("F6").Value = 1.068897904863723
$ws.Range("I6").Value = 1.043827855344991
$ws.Range("J6").Value = 1.061216202897955
$ws.Range("K6").Value = 1.062582678769481
$ws.Range("L6").Value = 1.055296515960157
$ws.Range("M6").Value = 1.071242991462381

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.05632475354435
$ws.Range("D7").Value = 1.059761118156985
$ws.Range("E7").Value = 1.052415559510553
$ws.Range("F7").Value = 1.068369156367585
$ws.Range("I7").Value = 1.043669344047017
$ws.Range("J7").Value = 1.060759102657733
$ws.Range("K7").Value = 1.062185768161252
$ws.Range("L7").Value = 1.05485802849697
$ws.Range("M7").Value = 1.070773259660677

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053942090866931
$ws.Range("D8").Value = 1.057847297775572
$ws.Range("E8").Value = 1.05032476109504
$ws.Range("F8").Value = 1.066152913502681
$ws.Range("I8").Value = 1.043000470553091
$ws.Range("J8").Value = 1.058840646259903
$ws.Range("K8").Value = 1.060519129628747
$ws.Range("L8").Value = 1.053016840751954
$ws.Range("M8").Value = 1.068802748432777

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.049714739621802
$ws.Range("D9").Value = 1.054449166956319
$ws.Range("E9").Value = 1.0466114657905
$ws.Range("F9").Value = 1.062223324188587
$ws.Range("I9").Value = 1.041797590977814
$ws.Range("J9").Value = 1.055429458626093
$ws.Range("K9").Value = 1.057552712841915
$ws.Range("L9").Value = 1.049739811779812
$ws.Range("M9").Value = 1.065302665906675

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04687595975614
$ws.Range("D10").Value = 1.052165562565986
$ws.Range("E10").Value = 1.044115371162573
$ws.Range("F10").Value = 1.059586261640353
$ws.Range("I10").Value = 1.040978987272023
$ws.Range("J10").Value = 1.053133775418172
$ws.Range("K10").Value = 1.055554380637495
$ws.Range("L10").Value = 1.047532222864535
$ws.Range("M10").Value = 1.062949653117611

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045641597991866
$ws.Range("D11").Value = 1.051172225885932
$ws.Range("E11").Value = 1.043029424309391
$ws.Range("F11").Value = 1.05844004650348
$ws.Range("I11").Value = 1.040620478489498
$ws.Range("J11").Value = 1.05213438857928
$ws.Range("K11").Value = 1.054683979983303
$ws.Range("L11").Value = 1.046570662750667
$ws.Range("M11").Value = 1.061925907054349

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.045182304129151
$ws.Range("D12").Value = 1.050802560020183
$ws.Range("E12").Value = 1.042625265688194
$ws.Range("F12").Value = 1.05801361807605
$ws.Range("I12").Value = 1.040486696597367
$ws.Range("J12").Value = 1.051762350543845
$ws.Range("K12").Value = 1.054359890600326
$ws.Range("L12").Value = 1.046212627364842
$ws.Range("M12").Value = 1.061544891313652

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045280860657812
$ws.Range("D13").Value = 1.050881886394825
$ws.Range("E13").Value = 1.042711995142481
$ws.Range("F13").Value = 1.05810511919754
$ws.Range("I13").Value = 1.040515421255382
$ws.Range("J13").Value = 1.051842191434824
$ws.Range("K13").Value = 1.054429444604671
$ws.Range("L13").Value = 1.046289466806598
$ws.Range("M13").Value = 1.061626654758395

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045603649009503
$ws.Range("D14").Value = 1.051141683496602
$ws.Range("E14").Value = 1.042996032657095
$ws.Range("F14").Value = 1.058404811641829
$ws.Range("I14").Value = 1.04060943264242
$ws.Range("J14").Value = 1.052103652662798
$ws.Range("K14").Value = 1.054657206737894
$ws.Range("L14").Value = 1.0465410852743
$ws.Range("M14").Value = 1.06189442759868

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045802423270015
$ws.Range("D15").Value = 1.051301660268453
$ws.Range("E15").Value = 1.043170932305516
$ws.Range("F15").Value = 1.058589372265891
$ws.Range("I15").Value = 1.040667274336176
$ws.Range("J15").Value = 1.052264638160422
$ws.Range("K15").Value = 1.054797434153516
$ws.Range("L15").Value = 1.04669599997792
$ws.Range("M15").Value = 1.062059311233693

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.046957769790819
$ws.Range("D16").Value = 1.052231390338559
$ws.Range("E16").Value = 1.044187332260473
$ws.Range("F16").Value = 1.059662238854432
$ws.Range("I16").Value = 1.041002694386311
$ws.Range("J16").Value = 1.053199987301171
$ws.Range("K16").Value = 1.055612037233868
$ws.Range("L16").Value = 1.047595917597793
$ws.Range("M16").Value = 1.063017491483491

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.047681092825488
$ws.Range("D17").Value = 1.052813363017271
$ws.Range("E17").Value = 1.044823507791791
$ws.Range("F17").Value = 1.060334041436654
$ws.Range("I17").Value = 1.041212005428691
$ws.Range("J17").Value = 1.053785263967518
$ws.Range("K17").Value = 1.056121636891837
$ws.Range("L17").Value = 1.048158883683357
$ws.Range("M17").Value = 1.063617214348591

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.048102499452583
$ws.Range("D18").Value = 1.053152382804491
$ws.Range("E18").Value = 1.045194085449995
$ws.Range("F18").Value = 1.060725474175843
$ws.Range("I18").Value = 1.041333702966753
$ws.Range("J18").Value = 1.054126131915769
$ws.Range("K18").Value = 1.056418385848874
$ws.Range("L18").Value = 1.048486707498536
$ws.Range("M18").Value = 1.063966553335063

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.048246104875314
$ws.Range("D19").Value = 1.053267906518503
$ws.Range("E19").Value = 1.045320359879789
$ws.Range("F19").Value = 1.060858872255759
$ws.Range("I19").Value = 1.041375132751188
$ws.Range("J19").Value = 1.054242272456048
$ws.Range("K19").Value = 1.056519486625203
$ws.Range("L19").Value = 1.048598395111163
$ws.Range("M19").Value = 1.064085589865091

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047603538523086
$ws.Range("D20").Value = 1.05275096797762
$ws.Range("E20").Value = 1.04475530329252
$ws.Range("F20").Value = 1.06026200678976
$ws.Range("I20").Value = 1.041189588734436
$ws.Range("J20").Value = 1.053722522624628
$ws.Range("K20").Value = 1.056067012667543
$ws.Range("L20").Value = 1.048098539235148
$ws.Range("M20").Value = 1.063552918425883

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045508618071049
$ws.Range("D21").Value = 1.051065199088407
$ws.Range("E21").Value = 1.042912412648847
$ws.Range("F21").Value = 1.058316578425941
$ws.Range("I21").Value = 1.040581765691919
$ws.Range("J21").Value = 1.052026681671574
$ws.Range("K21").Value = 1.05459015823478
$ws.Range("L21").Value = 1.046467014052367
$ws.Range("M21").Value = 1.061815596047211

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044186835949536
$ws.Range("D22").Value = 1.050001251478125
$ws.Range("E22").Value = 1.041749135545549
$ws.Range("F22").Value = 1.057089506206534
$ws.Range("I22").Value = 1.040196036733284
$ws.Range("J22").Value = 1.050955678042136
$ws.Range("K22").Value = 1.053657058099633
$ws.Range("L22").Value = 1.045436171283156
$ws.Range("M22").Value = 1.060718918962501

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04488798358006
$ws.Range("D23").Value = 1.050565658769393
$ws.Range("E23").Value = 1.042366251550456
$ws.Range("F23").Value = 1.057740377092088
$ws.Range("I23").Value = 1.04040085958123
$ws.Range("J23").Value = 1.051523894990701
$ws.Range("K23").Value = 1.054152148225783
$ws.Range("L23").Value = 1.045983124456914
$ws.Range("M23").Value = 1.061300707076699

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047638583504454
$ws.Range("D24").Value = 1.052779162954167
$ws.Range("E24").Value = 1.044786123492072
$ws.Range("F24").Value = 1.060294557438259
$ws.Range("I24").Value = 1.041199719078743
$ws.Range("J24").Value = 1.053750874324706
$ws.Range("K24").Value = 1.056091696522842
$ws.Range("L24").Value = 1.048125807976682
$ws.Range("M24").Value = 1.063581972435238

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050811144777749
$ws.Range("D25").Value = 1.05533080495763
$ws.Range("E25").Value = 1.047574989343325
$ws.Range("F25").Value = 1.063242198223796
$ws.Range("I25").Value = 1.042111476146531
$ws.Range("J25").Value = 1.056315060355958
$ws.Range("K25").Value = 1.058323193510335
$ws.Range("L25").Value = 1.050590967252723
$ws.Range("M25").Value = 1.066210910436378

